$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 13) following the same pattern as the
# existing SEO/accessibility rows above it. Values are assigned in this
# particular order so newly created shared-string entries line up with
# the order Excel itself produced them in (B, D, E, C).
$ws.Range("A13").Value = "Accessibilité"
$ws.Range("B13").Value = "Image a la place du texte"
$ws.Range("D13").Value = "Une image a la place d'un texte, pas responsive, pas optimisé"
$ws.Range("E13").Value = "Réécrire le texte "
$ws.Range("C13").Value = "Utiliser des paragraphes pour faire des citations"

# Update the view: scroll back to column A (remove the topLeftCell="C1"
# offset) and move the active selection to C26.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C26").Select()
